$wb = $excel.ActiveWorkbook

# --- Sheet "Full results" ---
$ws1 = $wb.Worksheets.Item("Full results")
$ws1.Range("K2").Value = 0.30548131530695
$ws1.Range("L2").Value = 0.0151543214118822
$ws1.Range("M2").Value = 0.0429552022729067
$ws1.Range("N2").Value = 0.320635636718833
$ws1.Range("F3").Value = 0.618516197488937
$ws1.Range("G3").Value = 0.305549957001496
$ws1.Range("H4").Value = 0.603358470898901
$ws1.Range("I4").Value = 0.283583444307945
$ws1.Range("O4").Value = 0.396777073377613
$ws1.Range("K5").Value = 0.118711403056588
$ws1.Range("L5").Value = 0.0284234699386578
$ws1.Range("M5").Value = 0.0518744755734492
$ws1.Range("N5").Value = 0.147134872995246
$ws1.Range("F6").Value = 0.853908072479899
$ws1.Range("G6").Value = 0.118720652051072
$ws1.Range("H7").Value = 0.825482388023428
$ws1.Range("I7").Value = 0.120961054165269
$ws1.Range("O7").Value = 0.174581921613705
$ws1.Range("K8").Value = 0.154950487430902
$ws1.Range("L8").Value = 0.0462753967250775
$ws1.Range("M8").Value = 0.0538964409094861
$ws1.Range("N8").Value = 0.20122588415598
$ws1.Range("F9").Value = 0.828052029768333
$ws1.Range("G9").Value = 0.154966664251263
$ws1.Range("H10").Value = 0.781771801894982
$ws1.Range("I10").Value = 0.130343263131413
$ws1.Range("O10").Value = 0.218309806508113
$ws1.Range("K11").Value = 0.227628738047098
$ws1.Range("L11").Value = 0.0295879891141118
$ws1.Range("M11").Value = 0.0375757067530702
$ws1.Range("N11").Value = 0.257216727161209
$ws1.Range("F12").Value = 0.760790786287031
$ws1.Range("G12").Value = 0.22766217657027
$ws1.Range("H13").Value = 0.731198450715314
$ws1.Range("I13").Value = 0.222597183806778
$ws1.Range("O13").Value = 0.268908946115844
$ws1.Range("K14").Value = 0.129157393728923
$ws1.Range("L14").Value = 0.0112827630153931
$ws1.Range("M14").Value = 0.0234822053534689
$ws1.Range("N14").Value = 0.140440156744316
$ws1.Range("F15").Value = 0.848362668967087
$ws1.Range("G15").Value = 0.129168831423602
$ws1.Range("H16").Value = 0.837078906792492
$ws1.Range("I16").Value = 0.135432813955002
$ws1.Range("O16").Value = 0.162995215213972

# --- Sheet "For plotting" ---
$ws2 = $wb.Worksheets.Item("For plotting")
$ws2.Range("B2").Value = 0.268541991328053
$ws2.Range("C2").Value = 0.372729282109613
$ws2.Range("D2").Value = 0.320635636718833
$ws2.Range("B3").Value = 0.350898215188114
$ws2.Range("C3").Value = 0.442655931567112
$ws2.Range("D3").Value = 0.396777073377613
$ws2.Range("B4").Value = 0.308550083518096
$ws2.Range("C4").Value = 0.399093658691316
$ws2.Range("B5").Value = 0.0917942259200481
$ws2.Range("C5").Value = 0.202475520070444
$ws2.Range("D5").Value = 0.147134872995246
$ws2.Range("B6").Value = 0.12393792036737
$ws2.Range("C6").Value = 0.225225922860039
$ws2.Range("D6").Value = 0.174581921613705
$ws2.Range("B7").Value = 0.0691131899280597
$ws2.Range("C7").Value = 0.176301702152451
$ws2.Range("B8").Value = 0.083920337627273
$ws2.Range("C8").Value = 0.318531430684686
$ws2.Range("D8").Value = 0.20122588415598
$ws2.Range("B9").Value = 0.106907985058182
$ws2.Range("C9").Value = 0.329711627958044
$ws2.Range("D9").Value = 0.218309806508113
$ws2.Range("B10").Value = 0.0385390742663888
$ws2.Range("C10").Value = 0.290287656930865
$ws2.Range("B11").Value = 0.120416345403155
$ws2.Range("C11").Value = 0.394017108919264
$ws2.Range("D11").Value = 0.257216727161209
$ws2.Range("B12").Value = 0.138419689548485
$ws2.Range("C12").Value = 0.399398202683203
$ws2.Range("D12").Value = 0.268908946115844
$ws2.Range("B13").Value = 0.0768065018400139
$ws2.Range("C13").Value = 0.385859976885533
$ws2.Range("B14").Value = 0.067715078846633
$ws2.Range("C14").Value = 0.213165234642
$ws2.Range("D14").Value = 0.140440156744316
$ws2.Range("B15").Value = 0.0925238746145066
$ws2.Range("C15").Value = 0.233466555813438
$ws2.Range("D15").Value = 0.162995215213972
$ws2.Range("B16").Value = 0.0663190443846013
$ws2.Range("C16").Value = 0.212706975336406
